# Generate Report for Handoff
# Updates the localization status report to reflect the "Ready for handoff"
# state (previously "In Translation"), bumps the recorded generate/handoff
# timestamps by one minute, and widens the Status columns to fit the new,
# longer status text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Per-language status (zh-cn / de-de) moves from "In Translation" to
# "Ready for handoff".
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"

# Latest HO Xliff Generate Date advances by a minute.
$wsOverview.Range("G2").Value = "2016-08-13 17:19:58"

# Widen the zh-cn / de-de status columns so the longer text fits.
$wsOverview.Range("E1").ColumnWidth = 16.33
$wsOverview.Range("F1").ColumnWidth = 16.33

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-13 17:19:50"
$wsZhCn.Range("C1").ColumnWidth = 16.33

# --- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-13 17:19:58"
$wsDeDe.Range("C1").ColumnWidth = 16.33
